# Apply updated "想去人数" (want-to-go count) and "最低票价" (min ticket price)
# figures to the "展览" (F1/F4 name conflict avoided), "演出" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 8700
$wsExpo.Range("G4").Value = 50
$wsExpo.Range("F9").Value = 475
$wsExpo.Range("F13").Value = 6254
$wsExpo.Range("F14").Value = 210
$wsExpo.Range("F15").Value = 324
$wsExpo.Range("F16").Value = 2391
$wsExpo.Range("F17").Value = 120
$wsExpo.Range("F18").Value = 205
$wsExpo.Range("F20").Value = 472

# --- Sheet "演出" (performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F5").Value = 47

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 8700
$wsAll.Range("G4").Value = 50
$wsAll.Range("F11").Value = 475
$wsAll.Range("F16").Value = 6254
$wsAll.Range("F17").Value = 47
$wsAll.Range("F18").Value = 210
$wsAll.Range("F19").Value = 324
$wsAll.Range("F20").Value = 2391
$wsAll.Range("F21").Value = 120
$wsAll.Range("F22").Value = 205
$wsAll.Range("F24").Value = 472
